$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cD = $ws.Cells.Item(2,4)
$cD.Value = "'29.140.67"
$cD.Style = "Normal"
$ws.Cells.Item(2,5).Value = "  -0.37%  "

# Row 3
$cD = $ws.Cells.Item(3,4)
$cD.Value = "'1.840.31"
$cD.Style = "Normal"
$ws.Cells.Item(3,5).Value = "  -0.50%  "

# Row 4
$cD = $ws.Cells.Item(4,4)
$cD.Value = "'0.9995"
$cD.Style = "Normal"
$ws.Cells.Item(4,5).Value = "  +0.04%  "

# Row 5
$cD = $ws.Cells.Item(5,4)
$cD.Value = "'241.16"
$cD.Style = "Normal"
$ws.Cells.Item(5,5).Value = "  -2.19%  "

# Row 6
$cD = $ws.Cells.Item(6,4)
$cD.Value = "'0.6867"
$cD.Style = "Normal"
$ws.Cells.Item(6,5).Value = "  -1.72%  "

# Row 7
$cD = $ws.Cells.Item(7,4)
$cD.Value = "'1.000"
$cD.Style = "Normal"
$ws.Cells.Item(7,5).Value = "  +0.05%  "

# Row 8
$ws.Cells.Item(8,5).Value = "  -1.46%  "

# Row 9
$cD = $ws.Cells.Item(9,4)
$cD.Value = "'0.07477"
$cD.Style = "Normal"
$ws.Cells.Item(9,5).Value = "  -3.23%  "

# Row 10
$cD = $ws.Cells.Item(10,4)
$cD.Value = "'23.18"
$cD.Style = "Normal"
$ws.Cells.Item(10,5).Value = "  -1.74%  "

# Row 11
$cD = $ws.Cells.Item(11,4)
$cD.Value = "'0.07665"
$cD.Style = "Normal"
$ws.Cells.Item(11,5).Value = "  -2.00%  "

# Row 12
$cD = $ws.Cells.Item(12,4)
$cD.Value = "'1.841.40"
$cD.Style = "Normal"
$ws.Cells.Item(12,5).Value = "  -0.36%  "

# Row 13
$cD = $ws.Cells.Item(13,4)
$cD.Value = "'5.068"
$cD.Style = "Normal"
$ws.Cells.Item(13,5).Value = "  -1.30%  "

# Row 14
$cD = $ws.Cells.Item(14,4)
$cD.Value = "'0.6836"
$cD.Style = "Normal"
$ws.Cells.Item(14,5).Value = "  -0.58%  "

# Row 15
$cD = $ws.Cells.Item(15,4)
$cD.Value = "'87.77"
$cD.Style = "Normal"
$ws.Cells.Item(15,5).Value = "  -6.05%  "

# Row 16
$cD = $ws.Cells.Item(16,4)
$cD.Value = "'6.163"
$cD.Style = "Normal"
$ws.Cells.Item(16,5).Value = "  -7.00%  "

# Row 17
$cD = $ws.Cells.Item(17,4)
$cD.Value = "'29.148.23"
$cD.Style = "Normal"
$ws.Cells.Item(17,5).Value = "  -0.27%  "

# Row 18
$cD = $ws.Cells.Item(18,4)
$cD.Value = "'0.000008171"
$cD.Style = "Normal"
$ws.Cells.Item(18,5).Value = "  -1.79%  "

# Row 19
$cD = $ws.Cells.Item(19,4)
$cD.Value = "'2.080.65"
$cD.Style = "Normal"
$ws.Cells.Item(19,5).Value = "  -0.45%  "

# Row 20
$cD = $ws.Cells.Item(20,4)
$cD.Value = "'228.02"
$cD.Style = "Normal"
$ws.Cells.Item(20,5).Value = "  -5.62%  "

# Row 21
$cD = $ws.Cells.Item(21,4)
$cD.Value = "'12.53"
$cD.Style = "Normal"
$ws.Cells.Item(21,5).Value = "  -1.96%  "

# Row 22
$cD = $ws.Cells.Item(22,4)
$cD.Value = "'0.9993"
$cD.Style = "Normal"
$ws.Cells.Item(22,5).Value = "  -0.04%  "

# Row 23
$cD = $ws.Cells.Item(23,4)
$cD.Value = "'7.416"
$cD.Style = "Normal"
$ws.Cells.Item(23,5).Value = "  -1.53%  "

# Row 24
$cD = $ws.Cells.Item(24,4)
$cD.Value = "'0.9999"
$cD.Style = "Normal"
$ws.Cells.Item(24,5).Value = "  +0.04%  "

# Row 25
$cD = $ws.Cells.Item(25,4)
$cD.Value = "'0.1459"
$cD.Style = "Normal"
$ws.Cells.Item(25,5).Value = "  -3.80%  "

# Row 26
$cD = $ws.Cells.Item(26,4)
$cD.Value = "'160.12"
$cD.Style = "Normal"

# Row 27
$cD = $ws.Cells.Item(27,4)
$cD.Value = "'8.773"
$cD.Style = "Normal"
$ws.Cells.Item(27,5).Value = "  -0.75%  "

# Row 28
$cD = $ws.Cells.Item(28,4)
$cD.Value = "'18.10"
$cD.Style = "Normal"
$ws.Cells.Item(28,5).Value = "  -1.04%  "

# Row 29
$cD = $ws.Cells.Item(29,4)
$cD.Value = "'1.516"
$cD.Style = "Normal"
$ws.Cells.Item(29,5).Value = "  -1.74%  "

# Row 30
$cD = $ws.Cells.Item(30,4)
$cD.Value = "'4.273"
$cD.Style = "Normal"
$ws.Cells.Item(30,5).Value = "  +0.81%  "

# Row 31
$cD = $ws.Cells.Item(31,4)
$cD.Value = "'4.156"
$cD.Style = "Normal"
$ws.Cells.Item(31,5).Value = "  -0.83%  "

# Row 32
$cD = $ws.Cells.Item(32,4)
$cD.Value = "'1.193"
$cD.Style = "Normal"
$ws.Cells.Item(32,5).Value = "  -0.41%  "

# Row 33
$cD = $ws.Cells.Item(33,4)
$cD.Value = "'0.05185"
$cD.Style = "Normal"
$ws.Cells.Item(33,5).Value = "  +1.24%  "

# Row 34
$ws.Cells.Item(34,5).Value = "  -3.18%  "

# Row 35
$cD = $ws.Cells.Item(35,4)
$cD.Value = "'1.854"
$cD.Style = "Normal"
$ws.Cells.Item(35,5).Value = "  -1.12%  "

# Row 36
$ws.Cells.Item(36,5).Value = "  -1.33%  "

# Row 37
$cD = $ws.Cells.Item(37,4)
$cD.Value = "'2.674"
$cD.Style = "Normal"
$ws.Cells.Item(37,5).Value = "  -0.61%  "

# Row 38
$cD = $ws.Cells.Item(38,4)
$cD.Value = "'1.317.07"
$cD.Style = "Normal"
$ws.Cells.Item(38,5).Value = "  +0.01%  "

# Row 39
$ws.Cells.Item(39,5).Value = "  -1.96%  "

# Row 40
$cD = $ws.Cells.Item(40,4)
$cD.Value = "'2.723"
$cD.Style = "Normal"
$ws.Cells.Item(40,5).Value = "  +0.48%  "

# Row 41
$cD = $ws.Cells.Item(41,4)
$cD.Value = "'0.9365"
$cD.Style = "Normal"
$ws.Cells.Item(41,5).Value = "  -1.69%  "

# Row 42
$cD = $ws.Cells.Item(42,4)
$cD.Value = "'104.78"
$cD.Style = "Normal"
$ws.Cells.Item(42,5).Value = "  -2.67%  "

# Row 43
$cD = $ws.Cells.Item(43,4)
$cD.Value = "'5.775"
$cD.Style = "Normal"
$ws.Cells.Item(43,5).Value = "  -4.53%  "

# Row 44
$cD = $ws.Cells.Item(44,4)
$cD.Value = "'0.9993"
$cD.Style = "Normal"
$ws.Cells.Item(44,5).Value = "  +0.00%  "

# Row 45
$cD = $ws.Cells.Item(45,4)
$cD.Value = "'0.00000000125"
$cD.Style = "Normal"
$ws.Cells.Item(45,5).Value = "  +1.19%  "

# Row 46
$ws.Cells.Item(46,2).Value = "RocketPoolETH"
$ws.Cells.Item(46,3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$cD = $ws.Cells.Item(46,4)
$cD.Value = "'1.982.84"
$cD.Style = "Normal"
$ws.Cells.Item(46,5).Value = "  -0.20%  "

# Row 47
$ws.Cells.Item(47,2).Value = "Aave"
$ws.Cells.Item(47,3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cD = $ws.Cells.Item(47,4)
$cD.Value = "'65.12"
$cD.Style = "Normal"
$ws.Cells.Item(47,5).Value = "  +1.33%  "

# Row 48
$cD = $ws.Cells.Item(48,4)
$cD.Value = "'0.5194"
$cD.Style = "Normal"
$ws.Cells.Item(48,5).Value = "  +0.22%  "

# Row 49
$cD = $ws.Cells.Item(49,4)
$cD.Value = "'9.545"
$cD.Style = "Normal"
$ws.Cells.Item(49,5).Value = "  -2.33%  "

# Row 50
$ws.Cells.Item(50,5).Value = "  +0.41%  "

# Row 51
$cD = $ws.Cells.Item(51,4)
$cD.Value = "'0.05944"
$cD.Style = "Normal"
$ws.Cells.Item(51,5).Value = "  +0.98%  "
